$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New daily rows appended to the "Chart" data table (dates, invalid-count,
# valid-count) - continuing directly after the existing 2025-12-02 row.
$newRows = @(
  @(60, "2025-12-03", 0, 26),
  @(61, "2025-12-04", 0, 25),
  @(62, "2025-12-05", 0, 25)
)

foreach ($row in $newRows) {
  $r = $row[0]
  $dateText = $row[1]
  $invalid = $row[2]
  $valid = $row[3]

  # Force the date-looking string to be stored as literal text (shared
  # string) instead of being auto-converted to a date serial number, then
  # drop the temporary Text number-format again so the cell ends up with
  # the same (default) style as all the other rows in the column.
  $dateCell = $ws.Cells.Item($r, 1)
  $dateCell.NumberFormat = "@"
  $dateCell.Value = $dateText
  $dateCell.ClearFormats()

  $ws.Cells.Item($r, 2).Value = $invalid
  $ws.Cells.Item($r, 3).Value = $valid
}
